$wb = $excel.ActiveWorkbook

# --- Friday sheet: remove the "Dates" row's search-snippet cells, move selection to C9 ---
$friday = $wb.Worksheets.Item("Friday")
$friday.Range("B2:C2").ClearContents()
$friday.Range("C9").Select()

# --- Saturday sheet: clear all the "Longest option" / "Shortest option" cells, move selection to B5 ---
$saturday = $wb.Worksheets.Item("Saturday")
$saturday.Range("B2:C11").ClearContents()
$saturday.Range("B5").Select()

# --- Sunday sheet: populate the "Longest option" / "Shortest option" cells, activate the sheet, move selection to B3 ---
$sunday = $wb.Worksheets.Item("Sunday")

$sunday.Range("B2").Value = "dates price in saudi arabia"
$sunday.Range("C2").Value = "Dates`nFruit"

$sunday.Range("B3").Value = "dhaka education board`nBoard of Intermediate and Secondary Education, Dhaka " + [char]0x00B7 + " 5 Joynag Rd, Dhaka"
$sunday.Range("C3").Value = "dhaka post"

$sunday.Range("B4").Value = "Baby Girl`nSong by Dhvani Bhanushali and Guru Randhawa"
$sunday.Range("C4").Value = "baby shark"

$sunday.Range("B5").Value = "School 2017`nTelevision series"
$sunday.Range("C5").Value = "school"

$sunday.Range("B6").Value = "cricket icc`nInternational Cricket Council " + [char]0x2014 + " Cricket administrative body"
$sunday.Range("C6").Value = "cricket"

$sunday.Range("B8").Value = "Inter Miami CF`nSoccer club"
$sunday.Range("C8").Value = "internet"

$sunday.Range("B9").Value = "look meaning in bengali"
$sunday.Range("C9").Value = "look"

$sunday.Range("B10").Value = "Hello Hello! (Noodle & Pals) [Sing-Along]`nSong by Noodle & Pals and Super Simple Songs"
$sunday.Range("C10").Value = "hello world"

$sunday.Range("B11").Value = "Byeon Woo-seok`nSouth Korean actor and model"
$sunday.Range("C11").Value = "by"

$sunday.Activate()
$sunday.Range("B3").Select()
